$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-15 down to 5-16
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with values
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 44530
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = 300000000
$ws.Cells.Item(4, 7).Value = "Espárragos"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 800
$ws.Cells.Item(4, 12).Value = 900
$ws.Cells.Item(4, 13).Value = 850
$ws.Cells.Item(4, 14).Value = "$/kilo"
$ws.Cells.Item(4, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(4, 16).Value = 850
$ws.Cells.Item(4, 17).Value = 1
$ws.Cells.Item(4, 18).Value = "Hortaliza"

# Copy the date style (s="2") from D5 (the row below, which carries on the previous D4 style) to the new D4
$ws.Cells.Item(5, 4).Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4122)  # xlPasteFormats
